# Scheduled market-data refresh for the Leve profit sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Updates currentAveragePrice / NQ / HQ price and
# profit columns (H:N) with freshly pulled values; a few rows lose their
# HQ-profit figure entirely (no HQ listings came back) and two rows on LTW
# gain a profit figure for the first time.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 863.4
$ws.Range("I18").Value = 329.25
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 329.25
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -45.25
$ws.Range("N18").Value = -3568
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H64").Value = 5213002.5
$ws.Range("I64").Value = 8336914.5
$ws.Range("J64").Value = 6482.5557
$ws.Range("K64").Value = 8336914.5
$ws.Range("L64").Value = 6482.5557
$ws.Range("M64").Value = -8336666.5
$ws.Range("N64").Value = -6978.5557
$ws.Range("H67").Value = 5213002.5
$ws.Range("I67").Value = 8336914.5
$ws.Range("J67").Value = 6482.5557
$ws.Range("K67").Value = 8336914.5
$ws.Range("L67").Value = 6482.5557
$ws.Range("M67").Value = -8336056.5
$ws.Range("N67").Value = -8198.555700000001
$ws.Range("H76").Value = 3475333
$ws.Range("I76").Value = 4447314.5
$ws.Range("J76").Value = 3971.4285
$ws.Range("K76").Value = 4447314.5
$ws.Range("L76").Value = 3971.4285
$ws.Range("M76").Value = -4446999.5
$ws.Range("N76").Value = -4601.4285
$ws.Range("H79").Value = 3475333
$ws.Range("I79").Value = 4447314.5
$ws.Range("J79").Value = 3971.4285
$ws.Range("K79").Value = 4447314.5
$ws.Range("L79").Value = 3971.4285
$ws.Range("M79").Value = -4446222.5
$ws.Range("N79").Value = -6155.4285
$ws.Range("H116").Value = 2304.5454
$ws.Range("I116").Value = 2799.1428
$ws.Range("K116").Value = 2799.1428
$ws.Range("M116").Value = 642.8571999999999
$ws.Range("H131").Value = 10888.046
$ws.Range("I131").Value = 3387.4614
$ws.Range("J131").Value = 21722.223
$ws.Range("K131").Value = 10162.3842
$ws.Range("L131").Value = 65166.66900000001
$ws.Range("M131").Value = -5122.3842
$ws.Range("N131").Value = -75246.66900000001
$ws.Range("H137").Value = 25002014
$ws.Range("I137").Value = 45455700
$ws.Range("K137").Value = 136367100
$ws.Range("M137").Value = -136364550
$ws.Range("H138").Value = 1916.26
$ws.Range("I138").Value = 902.6
$ws.Range("J138").Value = 2592.0334
$ws.Range("K138").Value = 2707.8
$ws.Range("L138").Value = 7776.100199999999
$ws.Range("M138").Value = 2432.2
$ws.Range("N138").Value = -18056.1002
$ws.Range("H141").Value = 4408.6665
$ws.Range("I141").Value = 2395.25
$ws.Range("J141").Value = 7337.273
$ws.Range("K141").Value = 7185.75
$ws.Range("L141").Value = 22011.819
$ws.Range("M141").Value = -2005.75
$ws.Range("N141").Value = -32371.819

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 18999.5
$ws.Range("I63").Value = 18999.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 18999.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -18313.5
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 18999.5
$ws.Range("I66").Value = 18999.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 94997.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -91565.5
$ws.Range("N66").ClearContents()
$ws.Range("H97").Value = 10752.1
$ws.Range("I97").Value = 14931.571
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 14931.571
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -14435.571
$ws.Range("N97").Value = -1992
$ws.Range("H109").Value = 45000
$ws.Range("J109").Value = 45000
$ws.Range("L109").Value = 45000
$ws.Range("N109").Value = -47774

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 676
$ws.Range("I25").Value = 676
$ws.Range("K25").Value = 676
$ws.Range("M25").Value = -441
$ws.Range("H86").Value = 1642.8572
$ws.Range("I86").Value = 1583.3334
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1583.3334
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -460.3334
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1642.8572
$ws.Range("I89").Value = 1583.3334
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 7916.666999999999
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -2300.666999999999
$ws.Range("N89").Value = -21232
$ws.Range("H105").Value = 3060.8298
$ws.Range("I105").Value = 2967.1875
$ws.Range("J105").Value = 3260.6
$ws.Range("K105").Value = 2967.1875
$ws.Range("L105").Value = 3260.6
$ws.Range("M105").Value = -1220.1875
$ws.Range("N105").Value = -6754.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3808.913
$ws.Range("I31").Value = 1727.1613
$ws.Range("J31").Value = 5507.184
$ws.Range("K31").Value = 1727.1613
$ws.Range("L31").Value = 5507.184
$ws.Range("M31").Value = -1432.1613
$ws.Range("N31").Value = -6097.184
$ws.Range("H34").Value = 3808.913
$ws.Range("I34").Value = 1727.1613
$ws.Range("J34").Value = 5507.184
$ws.Range("K34").Value = 1727.1613
$ws.Range("L34").Value = 5507.184
$ws.Range("M34").Value = -1525.1613
$ws.Range("N34").Value = -5911.184
$ws.Range("H58").Value = 37038636
$ws.Range("I58").Value = 41667864
$ws.Range("K58").Value = 41667864
$ws.Range("M58").Value = -41667661
$ws.Range("H62").Value = 26398
$ws.Range("I62").Value = 30810
$ws.Range("J62").Value = 8750
$ws.Range("K62").Value = 30810
$ws.Range("L62").Value = 8750
$ws.Range("M62").Value = -30186
$ws.Range("N62").Value = -9998
$ws.Range("H65").Value = 26398
$ws.Range("I65").Value = 30810
$ws.Range("J65").Value = 8750
$ws.Range("K65").Value = 154050
$ws.Range("L65").Value = 43750
$ws.Range("M65").Value = -150930
$ws.Range("N65").Value = -49990
$ws.Range("H132").Value = 4067089
$ws.Range("I132").Value = 5209790
$ws.Range("J132").Value = 4151.6665
$ws.Range("K132").Value = 15629370
$ws.Range("L132").Value = 12454.9995
$ws.Range("M132").Value = -15626840
$ws.Range("N132").Value = -17514.9995
$ws.Range("H134").Value = 44120716
$ws.Range("I134").Value = 66668668
$ws.Range("K134").Value = 200006004
$ws.Range("M134").Value = -200003469
$ws.Range("H136").Value = 37038636
$ws.Range("I136").Value = 41667864
$ws.Range("K136").Value = 125003592
$ws.Range("M136").Value = -125001042
$ws.Range("H139").Value = 86000
$ws.Range("J139").Value = 86000
$ws.Range("L139").Value = 86000
$ws.Range("N139").Value = -96280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 733631.0600000001
$ws.Range("I122").Value = 387.5
$ws.Range("J122").Value = 950888.4399999999
$ws.Range("K122").Value = 3487.5
$ws.Range("L122").Value = 8557995.959999999
$ws.Range("M122").Value = -1037.5
$ws.Range("N122").Value = -8562895.959999999
$ws.Range("H131").Value = 5465912.5
$ws.Range("I131").Value = 560
$ws.Range("J131").Value = 6537550.5
$ws.Range("K131").Value = 1680
$ws.Range("L131").Value = 19612651.5
$ws.Range("M131").Value = 3360
$ws.Range("N131").Value = -19622731.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4962.6
$ws.Range("I70").Value = 5007.073
$ws.Range("J70").Value = 4760
$ws.Range("K70").Value = 5007.073
$ws.Range("L70").Value = 4760
$ws.Range("M70").Value = -4737.073
$ws.Range("N70").Value = -5300
$ws.Range("H73").Value = 4962.6
$ws.Range("I73").Value = 5007.073
$ws.Range("J73").Value = 4760
$ws.Range("K73").Value = 5007.073
$ws.Range("L73").Value = 4760
$ws.Range("M73").Value = -4071.073
$ws.Range("N73").Value = -6632
$ws.Range("H80").Value = 43481150
$ws.Range("I80").Value = 2691.3845
$ws.Range("J80").Value = 100003150
$ws.Range("K80").Value = 2691.3845
$ws.Range("L80").Value = 100003150
$ws.Range("M80").Value = -1693.3845
$ws.Range("N80").Value = -100005146
$ws.Range("H83").Value = 43481150
$ws.Range("I83").Value = 2691.3845
$ws.Range("J83").Value = 100003150
$ws.Range("K83").Value = 13456.9225
$ws.Range("L83").Value = 500015750
$ws.Range("M83").Value = -8464.922500000001
$ws.Range("N83").Value = -500025734
$ws.Range("H132").Value = 2726.4443
$ws.Range("I132").Value = 2535.4583
$ws.Range("J132").Value = 3108.4167
$ws.Range("K132").Value = 7606.374899999999
$ws.Range("L132").Value = 9325.250100000001
$ws.Range("M132").Value = -5076.374899999999
$ws.Range("N132").Value = -14385.2501
$ws.Range("H139").Value = 49500
$ws.Range("J139").Value = 49500
$ws.Range("L139").Value = 49500
$ws.Range("N139").Value = -59780

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 1500
$ws.Range("K82").Value = 1500
$ws.Range("M82").Value = -1139
$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 1500
$ws.Range("K85").Value = 1500
$ws.Range("M85").Value = -252

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2713.5
$ws.Range("I132").Value = 2537.1707
$ws.Range("J132").Value = 3195.4666
$ws.Range("K132").Value = 7611.5121
$ws.Range("L132").Value = 9586.399800000001
$ws.Range("M132").Value = -5081.5121
$ws.Range("N132").Value = -14646.3998
$ws.Range("H136").Value = 3833.7073
$ws.Range("I136").Value = 1619.0435
$ws.Range("K136").Value = 4857.1305
$ws.Range("M136").Value = -2307.1305

Write-Host "Updated market price/profit columns across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."
